$d = $word.ActiveDocument

# The document currently has a single section with no header. Grab the
# default (primary) header of the last section -- wdHeaderFooterPrimary = 1.
$section = $d.Sections.Last
$header = $section.Headers.Item(1)

# Style + center the (already-present, empty) header paragraph.
$paragraph = $header.Range.Paragraphs(1)
$paragraph.Style = "Header"
$paragraph.Alignment = 1   # wdAlignParagraphCenter

# Insert the questionnaire-number text after the (empty) paragraph mark so
# the mark itself doesn't pick up direct character formatting.
$header.Range.InsertAfter("Questionnaire 45")

# Format only the inserted run (exclude the trailing paragraph mark).
$textRange = $header.Range.Duplicate
$textRange.End = $textRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
